# Cybersecurity Interface Agreement.docx - apply the "renamed product-level
# goals/requirements to element-level goals/requirements" edit:
#   1. Revision "Version 5" -> "Version 6"
#   2. Revision date stamp "4/22/24 12:16 PM" -> "4/11/25 2:37 PM"
#   3. Drop the stale lastRenderedPageBreak marker on the
#      "Agreement Harmonization" heading (gets cleared whenever the run
#      above it is touched / re-flowed)
#   4. "Product-level Security Requirements (AVCDL secondary document)" ->
#      "Element-level Security Requirements (AVCDL secondary document)"

$d = $word.ActiveDocument

# 1. Bump the revision number on the title page.
$r1 = $d.Content
$found = $r1.Find.Execute("Version 5", $true, $false, $false, $false, `
    $false, $true, 1, $false, "Version 6", 2)
if (-not $found) { throw "Could not find 'Version 5'" }

# 2. Update the cached result of the DATE field under the revision number
#    (leave the field code itself - " DATE \@ ""M/d/yy h:mm am/pm"" " -
#    untouched, just replace the cached/displayed text run).
$r2 = $d.Content
$found = $r2.Find.Execute("4/22/24 12:16 PM", $true, $false, $false, `
    $false, $false, $true, 1, $false, "4/11/25 2:37 PM", 2)
if (-not $found) { throw "Could not find the cached revision date" }

# 3. The "Agreement Harmonization" Heading2 run carried a stale
#    w:lastRenderedPageBreak from the previous layout pass; touching the
#    run (re-typing its text) drops that marker, matching the refreshed
#    pagination after the edits above.
$r3 = $d.Content
$found = $r3.Find.Execute("Agreement Harmonization", $true, $false, `
    $false, $false, $false, $true, 1, $false, "Agreement Harmonization", 2)
if (-not $found) { throw "Could not find 'Agreement Harmonization'" }

# 4. Rename "Product-level" to "Element-level" in the secondary-document
#    bullet to remove ambiguity with the AVCDL's own "product" terminology.
#    Only the leading "Product" word is retyped so the rest of the
#    sentence ("-level Security Requirements (AVCDL secondary document)")
#    is left untouched, same as the authored edit.
$r4 = $d.Content
$found = $r4.Find.Execute("Product-level Security Requirements (AVCDL secondary document)", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Could not find 'Product-level Security Requirements'" }
$start = $r4.Start
$wordRange = $d.Range($start, $start + 7)
if ($wordRange.Text -ne "Product") { throw "Unexpected text at target range: $($wordRange.Text)" }
$wordRange.Text = "Element"

Write-Output "Done"
